$wb = $excel.ActiveWorkbook

# Update the "nested" sheet text values: toBean -> toMap
$ws = $wb.Worksheets.Item("nested")
$ws.Range("C1").Value = "list#key?toMap=key"
$ws.Range("D1").Value = "list#value?toMap=value"

# Update selection on the "nested" sheet to C1
$ws.Activate()
$ws.Range("C1").Select()
